$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("preguntas")

# Add new "PISTA" header column (G1) and pista1..pista15 values (G2:G16)
$ws.Range("G1").Value = "PISTA"
for ($i = 1; $i -le 15; $i++) {
    $ws.Cells.Item($i + 1, 7).Value = "pista$i"
}

# Move the active selection to H8, matching the edited workbook's view state
$ws.Range("H8").Select()
